$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 7687.375
$ws.Range("I7").Value = 4499
$ws.Range("J7").Value = 8142.857
$ws.Range("K7").Value = 4499
$ws.Range("L7").Value = 8142.857
$ws.Range("M7").Value = -4387
$ws.Range("N7").Value = -8366.857

$ws.Range("H14").Value = 7687.375
$ws.Range("I14").Value = 4499
$ws.Range("J14").Value = 8142.857
$ws.Range("K14").Value = 4499
$ws.Range("L14").Value = 8142.857
$ws.Range("M14").Value = -4308
$ws.Range("N14").Value = -8524.857

$ws.Range("H19").Value = 893.5
$ws.Range("I19").Value = 359.7143
$ws.Range("J19").Value = 1071.4286
$ws.Range("K19").Value = 359.7143
$ws.Range("L19").Value = 1071.4286
$ws.Range("M19").Value = -184.7143
$ws.Range("N19").Value = -1421.4286

$ws.Range("H40").Value = 2485.3157
$ws.Range("I40").Value = 2538.4614
$ws.Range("K40").Value = 2538.4614
$ws.Range("M40").Value = -2363.4614

$ws.Range("H64").Value = 73919.36
$ws.Range("J64").Value = 2706.5454
$ws.Range("L64").Value = 2706.5454
$ws.Range("N64").Value = -3202.5454

$ws.Range("H67").Value = 73919.36
$ws.Range("J67").Value = 2706.5454
$ws.Range("L67").Value = 2706.5454
$ws.Range("N67").Value = -4422.5454

$ws.Range("H74").Value = 4375.769
$ws.Range("I74").Value = 5983.5713
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 5983.5713
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -5047.5713
$ws.Range("N74").Value = -4372

$ws.Range("H77").Value = 4375.769
$ws.Range("I77").Value = 5983.5713
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 29917.8565
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -25237.8565
$ws.Range("N77").Value = -21860

$ws.Range("H103").Value = 813.6875
$ws.Range("I103").Value = 2949.5
$ws.Range("J103").Value = 508.57144
$ws.Range("K103").Value = 8848.5
$ws.Range("L103").Value = 1525.71432
$ws.Range("M103").Value = -8262.5
$ws.Range("N103").Value = -2697.71432

$ws.Range("H129").Value = 2911.6345
$ws.Range("I129").Value = 11601.556
$ws.Range("J129").Value = 1092.814
$ws.Range("K129").Value = 34804.66800000001
$ws.Range("L129").Value = 3278.442
$ws.Range("M129").Value = -29804.66800000001
$ws.Range("N129").Value = -13278.442

$ws.Range("H137").Value = 1339.7826
$ws.Range("I137").Value = 1132.421
$ws.Range("K137").Value = 3397.263
$ws.Range("M137").Value = -847.2629999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 21199.6
$ws.Range("J9").Value = 21199.6
$ws.Range("L9").Value = 21199.6
$ws.Range("N9").Value = -21539.6

$ws.Range("H20").Value = 21199.6
$ws.Range("J20").Value = 21199.6
$ws.Range("L20").Value = 21199.6
$ws.Range("N20").Value = -21739.6

$ws.Range("H45").Value = 51691.2
$ws.Range("I45").Value = 63647.75
$ws.Range("K45").Value = 63647.75
$ws.Range("M45").Value = -63270.75

$ws.Range("H61").Value = 2287.9583
$ws.Range("I61").Value = 1446.4546
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1446.4546
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1234.4546
$ws.Range("N61").Value = -3424

$ws.Range("H63").Value = 2567.6667
$ws.Range("J63").Value = 2841.2
$ws.Range("L63").Value = 2841.2
$ws.Range("N63").Value = -4213.2

$ws.Range("H66").Value = 2567.6667
$ws.Range("J66").Value = 2841.2
$ws.Range("L66").Value = 14206
$ws.Range("N66").Value = -21070

$ws.Range("H102").Value = 68868.664
$ws.Range("I102").Value = 251872.5
$ws.Range("J102").Value = 2321.818
$ws.Range("K102").Value = 251872.5
$ws.Range("L102").Value = 2321.818
$ws.Range("M102").Value = -250250.5
$ws.Range("N102").Value = -5565.818

$ws.Range("H110").Value = 125250570
$ws.Range("I110").Value = 167000260
$ws.Range("J110").Value = 1506.5
$ws.Range("K110").Value = 167000260
$ws.Range("L110").Value = 1506.5
$ws.Range("M110").Value = -166998215
$ws.Range("N110").Value = -5596.5

$ws.Range("H122").Value = 1524.0358
$ws.Range("I122").Value = 1502.8077
$ws.Range("K122").Value = 4508.4231
$ws.Range("M122").Value = -2058.4231

$ws.Range("H136").Value = 2287.9583
$ws.Range("I136").Value = 1446.4546
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4339.3638
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1789.3638
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 688.25
$ws.Range("I94").Value = 592.375
$ws.Range("J94").Value = 880
$ws.Range("K94").Value = 592.375
$ws.Range("L94").Value = 880
$ws.Range("M94").Value = -141.375
$ws.Range("N94").Value = -1782

$ws.Range("H99").Value = 1553.8
$ws.Range("J99").Value = 2082.2632
$ws.Range("L99").Value = 2082.2632
$ws.Range("N99").Value = -5078.263199999999

$ws.Range("H105").Value = 97590.95
$ws.Range("I105").Value = 64880
$ws.Range("J105").Value = 202266
$ws.Range("K105").Value = 64880
$ws.Range("L105").Value = 202266
$ws.Range("M105").Value = -63133
$ws.Range("N105").Value = -205760

$ws.Range("H107").Value = 52676950
$ws.Range("I107").Value = 52676950
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 52676950
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -52675030
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value = 2938.4788
$ws.Range("I134").Value = 2677.5088
$ws.Range("J134").Value = 4001
$ws.Range("K134").Value = 8032.526400000001
$ws.Range("L134").Value = 12003
$ws.Range("M134").Value = -5497.526400000001
$ws.Range("N134").Value = -17073

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 909.36365
$ws.Range("I16").Value = 876.25
$ws.Range("K16").Value = 876.25
$ws.Range("M16").Value = -589.25

$ws.Range("H58").Value = 1112.8
$ws.Range("I58").Value = 1035.9025
$ws.Range("K58").Value = 1035.9025
$ws.Range("M58").Value = -832.9024999999999

$ws.Range("H99").Value = 10668.8125
$ws.Range("J99").Value = 16233.111
$ws.Range("L99").Value = 16233.111
$ws.Range("N99").Value = -19229.111

$ws.Range("H113").Value = 909.36365
$ws.Range("I113").Value = 876.25
$ws.Range("K113").Value = 876.25
$ws.Range("M113").Value = 1293.75

$ws.Range("H126").Value = 10668.8125
$ws.Range("J126").Value = 16233.111
$ws.Range("L126").Value = 48699.333
$ws.Range("N126").Value = -53639.333

$ws.Range("H136").Value = 1112.8
$ws.Range("I136").Value = 1035.9025
$ws.Range("K136").Value = 3107.7075
$ws.Range("M136").Value = -557.7074999999995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2430
$ws.Range("I82").Value = 1860
$ws.Range("K82").Value = 5580
$ws.Range("M82").Value = -5174

$ws.Range("H85").Value = 2430
$ws.Range("I85").Value = 1860
$ws.Range("K85").Value = 5580
$ws.Range("M85").Value = -4176

$ws.Range("H86").Value = 50000524
$ws.Range("I86").Value = 588
$ws.Range("J86").Value = 125000424
$ws.Range("K86").Value = 1764
$ws.Range("L86").Value = 375001272
$ws.Range("M86").Value = -578
$ws.Range("N86").Value = -375003644

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("N87").ClearContents()

$ws.Range("H89").Value = 50000524
$ws.Range("I89").Value = 588
$ws.Range("J89").Value = 125000424
$ws.Range("K89").Value = 5292
$ws.Range("L89").Value = 1125003816
$ws.Range("M89").Value = 636
$ws.Range("N89").Value = -1125015672

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("N90").ClearContents()

$ws.Range("H92").Value = 500
$ws.Range("J92").Value = 500
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 0
$ws.Range("N92").Value = -3996

$ws.Range("H131").Value = 6633.899
$ws.Range("J131").Value = 7197.4
$ws.Range("L131").Value = 21592.2
$ws.Range("N131").Value = -31672.2

$ws.Range("H140").Value = 6555.3335
$ws.Range("J140").Value = 1841.8889
$ws.Range("L140").Value = 5525.6667
$ws.Range("N140").Value = -15885.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 5250
$ws.Range("I25").Value = 3000
$ws.Range("J25").Value = 7500
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 7500
$ws.Range("M25").Value = -2471
$ws.Range("N25").Value = -8558

$ws.Range("H80").Value = 71504010
$ws.Range("I80").Value = 125130376
$ws.Range("J80").Value = 2191.6667
$ws.Range("K80").Value = 125130376
$ws.Range("L80").Value = 2191.6667
$ws.Range("M80").Value = -125129378
$ws.Range("N80").Value = -4187.6667

$ws.Range("H83").Value = 71504010
$ws.Range("I83").Value = 125130376
$ws.Range("J83").Value = 2191.6667
$ws.Range("K83").Value = 625651880
$ws.Range("L83").Value = 10958.3335
$ws.Range("M83").Value = -625646888
$ws.Range("N83").Value = -20942.3335

$ws.Range("H102").Value = 2617.353
$ws.Range("I102").Value = 1850.125
$ws.Range("J102").Value = 3299.3333
$ws.Range("K102").Value = 1850.125
$ws.Range("L102").Value = 3299.3333
$ws.Range("M102").Value = -228.125
$ws.Range("N102").Value = -6543.3333

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -7900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3231.6
$ws.Range("I7").Value = 1889.4
$ws.Range("J7").Value = 5916
$ws.Range("K7").Value = 1889.4
$ws.Range("L7").Value = 5916
$ws.Range("M7").Value = -1777.4
$ws.Range("N7").Value = -6140

$ws.Range("H40").Value = 60958.53
$ws.Range("I40").Value = 1000000
$ws.Range("J40").Value = 2268.4375
$ws.Range("K40").Value = 1000000
$ws.Range("L40").Value = 2268.4375
$ws.Range("M40").Value = -999864
$ws.Range("N40").Value = -2540.4375

$ws.Range("H61").Value = 2218
$ws.Range("I61").Value = 2525.75
$ws.Range("J61").Value = 1602.5
$ws.Range("K61").Value = 2525.75
$ws.Range("L61").Value = 1602.5
$ws.Range("M61").Value = -2323.75
$ws.Range("N61").Value = -2006.5

$ws.Range("H68").Value = 2893.9
$ws.Range("I68").Value = 1565.3334
$ws.Range("J68").Value = 3463.2856
$ws.Range("K68").Value = 1565.3334
$ws.Range("L68").Value = 3463.2856
$ws.Range("M68").Value = -816.3334
$ws.Range("N68").Value = -4961.2856

$ws.Range("H71").Value = 2893.9
$ws.Range("I71").Value = 1565.3334
$ws.Range("J71").Value = 3463.2856
$ws.Range("K71").Value = 7826.666999999999
$ws.Range("L71").Value = 17316.428
$ws.Range("M71").Value = -4082.666999999999
$ws.Range("N71").Value = -24804.428

$ws.Range("H82").Value = 1630.2941
$ws.Range("I82").Value = 857.1429000000001
$ws.Range("J82").Value = 2171.5
$ws.Range("K82").Value = 857.1429000000001
$ws.Range("L82").Value = 2171.5
$ws.Range("M82").Value = -496.1429000000001
$ws.Range("N82").Value = -2893.5

$ws.Range("H85").Value = 1630.2941
$ws.Range("I85").Value = 857.1429000000001
$ws.Range("J85").Value = 2171.5
$ws.Range("K85").Value = 857.1429000000001
$ws.Range("L85").Value = 2171.5
$ws.Range("M85").Value = 390.8570999999999
$ws.Range("N85").Value = -4667.5

$ws.Range("H100").Value = 2377.1333
$ws.Range("I100").Value = 2142.818
$ws.Range("J100").Value = 3021.5
$ws.Range("K100").Value = 2142.818
$ws.Range("L100").Value = 3021.5
$ws.Range("M100").Value = -1601.818
$ws.Range("N100").Value = -4103.5

$ws.Range("H113").Value = 2218
$ws.Range("I113").Value = 2525.75
$ws.Range("J113").Value = 1602.5
$ws.Range("K113").Value = 2525.75
$ws.Range("L113").Value = 1602.5
$ws.Range("M113").Value = -355.75
$ws.Range("N113").Value = -5942.5

$ws.Range("H126").Value = 3231.6
$ws.Range("I126").Value = 1889.4
$ws.Range("J126").Value = 5916
$ws.Range("K126").Value = 5668.200000000001
$ws.Range("L126").Value = 17748
$ws.Range("M126").Value = -3198.200000000001
$ws.Range("N126").Value = -22688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 18000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 18000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 18000
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -18346

$ws.Range("H96").Value = 62501490
$ws.Range("I96").Value = 142858720
$ws.Range("J96").Value = 1410.4445
$ws.Range("K96").Value = 142858720
$ws.Range("L96").Value = 1410.4445
$ws.Range("M96").Value = -142857347
$ws.Range("N96").Value = -4156.4445

$ws.Range("H122").Value = 1848.1428
$ws.Range("I122").Value = 1602.6666
$ws.Range("J122").Value = 1889.0555
$ws.Range("K122").Value = 4807.9998
$ws.Range("L122").Value = 5667.166499999999
$ws.Range("M122").Value = -2357.9998
$ws.Range("N122").Value = -10567.1665

$ws.Range("H126").Value = 1960.5555
$ws.Range("I126").Value = 1852
$ws.Range("J126").Value = 2177.6667
$ws.Range("K126").Value = 5556
$ws.Range("L126").Value = 6533.000100000001
$ws.Range("M126").Value = -3086
$ws.Range("N126").Value = -11473.0001

$ws.Range("H136").Value = 608.0925999999999
$ws.Range("I136").Value = 367.325
$ws.Range("J136").Value = 1296
$ws.Range("K136").Value = 1101.975
$ws.Range("L136").Value = 3888
$ws.Range("M136").Value = 1448.025
$ws.Range("N136").Value = -8988
